# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Re-order the merged cell ranges (A47:D47 now listed before A1:D1)
$ws.Range("A1:D1").UnMerge()
$ws.Range("A47:D47").UnMerge()
$ws.Range("A47:D47").Merge()
$ws.Range("A1:D1").Merge()

# Update the date in A1 (45308 -> 45309, i.e. 2024-01-17 -> 2024-01-18)
$ws.Range("A1").Value = 45309

# --- CAÑOS EXTENSIBLES section (rows 10-17) ---
$ws.Range("D10").Value = 3231
$ws.Range("D11").Value = 3913
$ws.Range("D12").Value = 3640
$ws.Range("D13").Value = 4459
$ws.Range("D14").Value = 5187
$ws.Range("D15").Value = 3822
$ws.Range("D16").Value = 6753
$ws.Range("D17").Value = 324

# --- CAÑOS CURVOS (Barrales) section (rows 25-32) ---
$ws.Range("D25").Value = 4050
$ws.Range("D26").Value = 7080
$ws.Range("D27").Value = 5460
$ws.Range("D28").Value = 7470
$ws.Range("D29").Value = 5100
$ws.Range("D30").Value = 8477
$ws.Range("D31").Value = 6734
$ws.Range("D32").Value = 8750
